$d = $word.ActiveDocument

# 1) Title-case "Mayornan/Dunadónan di Kuido" -> "Mayornan/Edukadónan"
#    (covers the document title and the two section headings)
$d.Content.Find.Execute("Mayornan/Dunadónan di Kuido", $true, $false, $false, $false, $false, $true, 1, $false, "Mayornan/Edukadónan", 2)

# 2) One special mixed-case occurrence: "Mayornan/dunadónan di kuido" -> "Mayornan/Edukadónan"
#    (capital M, lowercase d in source, but target capitalizes the E)
$d.Content.Find.Execute("Mayornan/dunadónan di kuido", $true, $false, $false, $false, $false, $true, 1, $false, "Mayornan/Edukadónan", 2)

# 3) Singular "mayor/dunadónan di kuido" -> "mayor/edukadónan"
$d.Content.Find.Execute("mayor/dunadónan di kuido", $true, $false, $false, $false, $false, $true, 1, $false, "mayor/edukadónan", 2)

# 4) Lowercase plural "mayornan/dunadónan di kuido" -> "mayornan/edukadónan" (bulk of occurrences)
$d.Content.Find.Execute("mayornan/dunadónan di kuido", $true, $false, $false, $false, $false, $true, 1, $false, "mayornan/edukadónan", 2)

# 5) Typo fix: "yunann" -> "yunan"
$d.Content.Find.Execute("yunann", $true, $false, $false, $false, $false, $true, 1, $false, "yunan", 2)

# 6) Remove the space in "Pasa Tempu hasiendo Hòbi" -> "PasaTempu hasiendo Hòbi"
$d.Content.Find.Execute("Pasa Tempu hasiendo Hòbi", $true, $false, $false, $false, $false, $true, 1, $false, "PasaTempu hasiendo Hòbi", 2)
